# Apply the update described by the commit:
# "updates made to application_data_condensed_df"
#
# A new column, WEEKEND_APPR_PROCESS_START, is inserted right before the
# existing FRAUD_RISK column. For each data row its value is the weekday
# name (taken from WEEKDAY_APPR_PROCESS_START) when that day falls on a
# weekend (SATURDAY/SUNDAY), otherwise "N/A".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the FRAUD_RISK column dynamically by scanning row 1.
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count
$lastRow = $usedRange.Rows.Count

$fraudCol = 0
$weekdayCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $hdr = $ws.Cells.Item(1, $c).Value()
    if ($hdr -eq "FRAUD_RISK") {
        $fraudCol = $c
    }
    if ($hdr -eq "WEEKDAY_APPR_PROCESS_START") {
        $weekdayCol = $c
    }
}

# Insert a new blank column immediately before FRAUD_RISK; this shifts
# FRAUD_RISK (and anything after it) one column to the right.
$ws.Columns.Item($fraudCol).Insert()

$newCol = $fraudCol

# Header for the newly inserted column. Columns.Insert() already carries
# over the formatting (border/bold/alignment) from the displaced column,
# so only the text needs to be set.
$ws.Cells.Item(1, $newCol).Value() = "WEEKEND_APPR_PROCESS_START"

# Fill data rows based on the weekday column.
for ($r = 2; $r -le $lastRow; $r++) {
    $weekday = $ws.Cells.Item($r, $weekdayCol).Value()
    if ($weekday -eq "SATURDAY" -or $weekday -eq "SUNDAY") {
        $ws.Cells.Item($r, $newCol).Value() = $weekday
    } else {
        $ws.Cells.Item($r, $newCol).Value() = "N/A"
    }
}
